$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.047456
$ws.Range("H2").Value = 0.142368
$ws.Range("I2").Value = 0.03340144944197188
$ws.Range("J2").Value = 0.03340144944197188
$ws.Range("M2").Value = 0.655792
$ws.Range("N2").Value = 1.967376
$ws.Range("O2").Value = 0.01246532615150124
$ws.Range("P2").Value = 0.01246532615150124
$ws.Range("Q2").Value = 0.031121265152
$ws.Range("R2").Value = 0.280091386368
$ws.Range("S2").Value = 0.0004163599612270588
$ws.Range("T2").Value = 0.0004163599612270587
# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.047456
$ws.Range("H3").Value = 0.142368
$ws.Range("I3").Value = 0.03340144944197188
$ws.Range("J3").Value = 0.03340144944197188
$ws.Range("O3").Value = 0.5315769812025607
$ws.Range("P3").Value = 0.5315769812025607
$ws.Range("Q3").Value = 1.327149244202667
$ws.Range("R3").Value = 11.944343197824
$ws.Range("S3").Value = 0.01775544166215337
$ws.Range("T3").Value = 0.01775544166215336
# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.047456
$ws.Range("H4").Value = 0.142368
$ws.Range("I4").Value = 0.03340144944197188
$ws.Range("J4").Value = 0.03340144944197188
$ws.Range("M4").Value = 23.987612
$ws.Range("N4").Value = 71.962836
$ws.Range("O4").Value = 0.4559576926459381
$ws.Range("P4").Value = 0.4559576926459381
$ws.Range("Q4").Value = 1.138356115072
$ws.Range("R4").Value = 10.245205035648
$ws.Range("S4").Value = 0.01522964781859145
$ws.Range("T4").Value = 0.01522964781859145
# Row 5
$ws.Range("I5").Value = 0.7427665150281654
$ws.Range("J5").Value = 0.7427665150281653
$ws.Range("M5").Value = 0.655792
$ws.Range("N5").Value = 1.967376
$ws.Range("O5").Value = 0.01246532615150124
$ws.Range("P5").Value = 0.01246532615150124
$ws.Range("Q5").Value = 0.6920607951573333
$ws.Range("R5").Value = 6.228547156416
$ws.Range("S5").Value = 0.009258826864240031
$ws.Range("T5").Value = 0.00925882686424003
# Row 6
$ws.Range("I6").Value = 0.7427665150281654
$ws.Range("J6").Value = 0.7427665150281653
$ws.Range("O6").Value = 0.5315769812025607
$ws.Range("P6").Value = 0.5315769812025607
$ws.Range("S6").Value = 0.3948375817970186
$ws.Range("T6").Value = 0.3948375817970185
# Row 7
$ws.Range("I7").Value = 0.7427665150281654
$ws.Range("J7").Value = 0.7427665150281653
$ws.Range("M7").Value = 23.987612
$ws.Range("N7").Value = 71.962836
$ws.Range("O7").Value = 0.4559576926459381
$ws.Range("P7").Value = 0.4559576926459381
$ws.Range("Q7").Value = 25.31425487753067
$ws.Range("R7").Value = 227.828293897776
$ws.Range("S7").Value = 0.3386701063669068
$ws.Range("T7").Value = 0.3386701063669068
# Row 8
$ws.Range("G8").Value = 0.3180153333333334
$ws.Range("H8").Value = 0.9540460000000001
$ws.Range("I8").Value = 0.2238320355298628
$ws.Range("J8").Value = 0.2238320355298628
$ws.Range("M8").Value = 0.655792
$ws.Range("N8").Value = 1.967376
$ws.Range("O8").Value = 0.01246532615150124
$ws.Range("P8").Value = 0.01246532615150124
$ws.Range("Q8").Value = 0.2085519114773334
$ws.Range("R8").Value = 1.876967203296
$ws.Range("S8").Value = 0.002790139326034155
$ws.Range("T8").Value = 0.002790139326034154
# Row 9
$ws.Range("G9").Value = 0.3180153333333334
$ws.Range("H9").Value = 0.9540460000000001
$ws.Range("I9").Value = 0.2238320355298628
$ws.Range("J9").Value = 0.2238320355298628
$ws.Range("O9").Value = 0.5315769812025607
$ws.Range("P9").Value = 0.5315769812025607
$ws.Range("Q9").Value = 8.893581618303113
$ws.Range("R9").Value = 80.04223456472801
$ws.Range("S9").Value = 0.1189839577433888
$ws.Range("T9").Value = 0.1189839577433888
# Row 10
$ws.Range("G10").Value = 0.3180153333333334
$ws.Range("H10").Value = 0.9540460000000001
$ws.Range("I10").Value = 0.2238320355298628
$ws.Range("J10").Value = 0.2238320355298628
$ws.Range("M10").Value = 23.987612
$ws.Range("N10").Value = 71.962836
$ws.Range("O10").Value = 0.4559576926459381
$ws.Range("P10").Value = 0.4559576926459381
$ws.Range("Q10").Value = 7.628428426050667
$ws.Range("R10").Value = 68.65585583445601
$ws.Range("S10").Value = 0.1020579384604399
$ws.Range("T10").Value = 0.1020579384604399
